$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 6) mirroring the layout of the existing rows.
$row = 6

$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value = "Bíobío"

# Column D carries the same date style as the rows above (numFmt 165).
$ws.Cells.Item($row, 4).Value = 45106
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100104
$ws.Cells.Item($row, 8).Value = "Frutos de pepita"
$ws.Cells.Item($row, 9).Value = 100104001
$ws.Cells.Item($row, 10).Value = "Granada"
$ws.Cells.Item($row, 11).Value = "Wonderfull"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 50
$ws.Cells.Item($row, 14).Value = 10000
$ws.Cells.Item($row, 15).Value = 10000
$ws.Cells.Item($row, 16).Value = 10000
$ws.Cells.Item($row, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($row, 19).Value = 556
$ws.Cells.Item($row, 20).Value = 18
